# PV.xlsx - Include PV function description
# Adds a "Formula Text" column (H) that shows the formula used in column F,
# and a new error-propagation example table (rows 9-14) with a
# Formula / Formula Text / Comments layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Formula Text" column next to the existing table (H1:H6)
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Formula Text"

$ws.Range("H2").Formula = "=FORMULATEXT(F2)"
$ws.Range("H3:H6").Formula = "=FORMULATEXT(F3)"

# ---------------------------------------------------------------------
# 2. New error-propagation example table (rows 9-14)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Formula"
$ws.Range("B9").Value = "Formula Text"
$ws.Range("D9").Value = "Comments"
$ws.Range("A9:E9").HorizontalAlignment = -4108
$ws.Range("A9:E9").Font.Bold = $true

# -- Row 10: SQRT(-1) -> #NUM!
$ws.Range("A10").Formula = "=PV(SQRT(-1),2,3)"
$ws.Range("A10").NumberFormat = '"£"#,##0.00;[Red]\-"£"#,##0.00'
$ws.Range("B10").Formula = "=FORMULATEXT(A10)"
$ws.Range("B10").HorizontalAlignment = -4131

# -- Row 11: text instead of rate -> #VALUE!
$ws.Range("A11").Formula = '=PV("str",2,3,4,TRUE)'
$ws.Range("A11").NumberFormat = '"£"#,##0.00;[Red]\-"£"#,##0.00'
$ws.Range("B11").Formula = "=FORMULATEXT(A11)"
$ws.Range("B11").HorizontalAlignment = -4131

# -- Row 12: text instead of type -> #VALUE! (shared formula group starts)
$ws.Range("A12").Formula = '=PV(1,2,3,4,"false")'
$ws.Range("A12").NumberFormat = '"£"#,##0.00;[Red]\-"£"#,##0.00'
$ws.Range("B12:B14").Formula = "=FORMULATEXT(A12)"
$ws.Range("B12:B14").HorizontalAlignment = -4131

# -- Row 13: combination causing #NUM!
$ws.Range("A13").Formula = '=PV(-3,0.5,1)'
$ws.Range("A13").NumberFormat = '"£"#,##0.00;[Red]\-"£"#,##0.00'

# -- Row 14: combination causing #DIV/0!
$ws.Range("A14").Formula = '=PV(-1,20,300)'
$ws.Range("A14").NumberFormat = '"£"#,##0.00;[Red]\-"£"#,##0.00'

# -- Comments column (D:J merged on each row), wrapped text
$ws.Range("D10").Value = "Example of error propagation."
$ws.Range("D11").Value = "Unable to convert rate argument to a number."
$ws.Range("D11").Characters(1, 18).Font.Italic = $false
$ws.Range("D11").Characters(19, 4).Font.Italic = $true
$ws.Range("D11").Characters(23, 23).Font.Italic = $false
$ws.Range("D12").Value = "Unable to convert type argument to a Boolean."
$ws.Range("D12").Characters(1, 18).Font.Italic = $false
$ws.Range("D12").Characters(19, 4).Font.Italic = $true
$ws.Range("D12").Characters(23, 23).Font.Italic = $false
$ws.Range("D13").Value = "Combination of valid inputs cause a #NUM! error. In this case, the formula for FV includes the square root of -2."
$ws.Range("D14").Value = "Combination of valid inputs cause a #DIV/0! error. In this case, the formula for FV includes 1 divided by 0."

$ws.Range("D10:J14").WrapText = $true

$ws.Range("D10:J10").Merge()
$ws.Range("D11:J11").Merge()
$ws.Range("D12:J12").Merge()
$ws.Range("D13:J13").Merge()
$ws.Range("D14:J14").Merge()

# ---------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 10.15
$ws.Columns("H").ColumnWidth = 26.15
$ws.Columns("J").ColumnWidth = 14.5

# ---------------------------------------------------------------------
# 4. Selection
# ---------------------------------------------------------------------
$ws.Range("A11").Select()
